$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 801, shifting existing rows 801-842 down to 802-843.
$ws.Rows("801:801").Insert()

# Populate the newly inserted row 801 with the new record.
# Temporarily force column A to text format so the date-like string
# "2026/02/12" is stored as literal text rather than being auto-converted
# to a date serial number, then clear the format back to the sheet's
# default (unstyled) look.
$ws.Range("A801").NumberFormat = "@"
$ws.Range("A801").Value = "2026/02/12"
$ws.Range("A801").ClearFormats()

$ws.Range("B801").Value = "木"
$ws.Range("C801").Value = 18
$ws.Range("D801").Value = 201
